$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 0.2857142857142857
$ws.Cells.Item(4, 4).Value = 0.25
$ws.Cells.Item(5, 4).Value = 0.3571428571428572
$ws.Cells.Item(6, 4).Value = 0.3214285714285715
$ws.Cells.Item(7, 4).Value = 0.4285714285714286
$ws.Cells.Item(8, 4).Value = 0.3839285714285715
$ws.Cells.Item(9, 4).Value = 0.2857142857142858
$ws.Cells.Item(10, 4).Value = 0.3169642857142858
$ws.Cells.Item(11, 4).Value = 0.3660714285714286
$ws.Cells.Item(12, 4).Value = 0.1428571428571428
$ws.Cells.Item(13, 4).Value = 0.07142857142857142
$ws.Cells.Item(14, 4).Value = 0.3035714285714286
$ws.Cells.Item(15, 4).Value = 0.4464285714285715
$ws.Cells.Item(16, 4).Value = 0.25
$ws.Cells.Item(17, 4).Value = 0.375
$ws.Cells.Item(18, 4).Value = 0.3526785714285714
$ws.Cells.Item(19, 4).Value = 0.3705357142857143
$ws.Cells.Item(20, 4).Value = 0.3236607142857143
$ws.Cells.Item(21, 4).Value = 0.3861607142857144
$ws.Cells.Item(22, 4).Value = 0.2857142857142857
$ws.Cells.Item(23, 4).Value = 0.3928571428571428
$ws.Cells.Item(24, 4).Value = 0.1428571428571428
$ws.Cells.Item(25, 4).Value = 0.125
$ws.Cells.Item(26, 4).Value = 0.4732142857142857
$ws.Cells.Item(27, 4).Value = 0.3214285714285715
$ws.Cells.Item(28, 4).Value = 0.3080357142857142
$ws.Cells.Item(29, 4).Value = 0.3616071428571428
$ws.Cells.Item(30, 4).Value = 0.3571428571428571
$ws.Cells.Item(31, 4).Value = 0.3705357142857142
$ws.Cells.Item(32, 4).Value = 0.1785714285714286
$ws.Cells.Item(33, 4).Value = 0.3214285714285714
$ws.Cells.Item(34, 4).Value = 0.3035714285714285
$ws.Cells.Item(35, 4).Value = 0.4285714285714286
$ws.Cells.Item(36, 4).Value = 0.25
$ws.Cells.Item(37, 4).Value = 0.2767857142857142
$ws.Cells.Item(38, 4).Value = 0.4017857142857144
$ws.Cells.Item(39, 4).Value = 0.2589285714285714
$ws.Cells.Item(40, 4).Value = 0.3459821428571428
$ws.Cells.Item(41, 4).Value = 0.3727678571428571
$ws.Cells.Item(42, 4).Value = 0.2142857142857143
$ws.Cells.Item(43, 4).Value = 0.2142857142857143
$ws.Cells.Item(44, 4).Value = 0.3035714285714286
$ws.Cells.Item(45, 4).Value = 0.2678571428571428
$ws.Cells.Item(46, 4).Value = 0.3303571428571428
$ws.Cells.Item(47, 4).Value = 0.5178571428571428
$ws.Cells.Item(48, 4).Value = 0.4062499999999999
$ws.Cells.Item(49, 4).Value = 0.4464285714285714
$ws.Cells.Item(50, 4).Value = 0.2946428571428572
$ws.Cells.Item(51, 4).Value = 0.3593749999999999
$ws.Cells.Item(52, 4).Value = 0.1785714285714286
$ws.Cells.Item(53, 4).Value = 0.07142857142857142
$ws.Cells.Item(54, 4).Value = 0.3392857142857143
$ws.Cells.Item(55, 4).Value = 0.3571428571428572
$ws.Cells.Item(56, 4).Value = 0.3571428571428572
$ws.Cells.Item(57, 4).Value = 0.3392857142857142
$ws.Cells.Item(58, 4).Value = 0.3794642857142856
$ws.Cells.Item(59, 4).Value = 0.3883928571428572
$ws.Cells.Item(60, 4).Value = 0.3526785714285714
$ws.Cells.Item(61, 4).Value = 0.3973214285714285
